# Update the multiplication-table answers to the new values.
# The document contains one table with 20 rows; every 5th row (1, 5, 10,
# 15, 20) holds five answer cells, the rest are spacer rows.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1 (table row 1)
$t.Cell(1, 1).Range.Text = "51×84=4284"
$t.Cell(1, 2).Range.Text = "51×45=2295"
$t.Cell(1, 3).Range.Text = "53×90=4770"
$t.Cell(1, 4).Range.Text = "57×11=627"
$t.Cell(1, 5).Range.Text = "81×84=6804"

# Row 2 (table row 5)
$t.Cell(5, 1).Range.Text = "26×97=2522"
$t.Cell(5, 2).Range.Text = "96×75=7200"
$t.Cell(5, 3).Range.Text = "24×57=1368"
$t.Cell(5, 4).Range.Text = "84×48=4032"
$t.Cell(5, 5).Range.Text = "21×35=735"

# Row 3 (table row 10)
$t.Cell(10, 1).Range.Text = "54×15=810"
$t.Cell(10, 2).Range.Text = "39×49=1911"
$t.Cell(10, 3).Range.Text = "20×65=1300"
$t.Cell(10, 4).Range.Text = "91×97=8827"
$t.Cell(10, 5).Range.Text = "30×37=1110"

# Row 4 (table row 15)
$t.Cell(15, 1).Range.Text = "89×80=7120"
$t.Cell(15, 2).Range.Text = "50×12=600"
$t.Cell(15, 3).Range.Text = "54×96=5184"
$t.Cell(15, 4).Range.Text = "35×51=1785"
$t.Cell(15, 5).Range.Text = "15×87=1305"

# Row 5 (table row 20)
$t.Cell(20, 1).Range.Text = "17×45=765"
$t.Cell(20, 2).Range.Text = "53×47=2491"
$t.Cell(20, 3).Range.Text = "45×97=4365"
$t.Cell(20, 4).Range.Text = "27×27=729"
$t.Cell(20, 5).Range.Text = "87×42=3654"
